$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking") changes
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total") changes
$ws.Range("B12").Value = 40
$ws.Range("C12").Value = -8
$ws.Range("E12").Value = "32 / 112"
